$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 85, pushing existing rows 85:160 down to 86:161.
$ws.Rows("85:85").Insert()

# The row that used to be row 85 is now row 86; copy its "constant" columns
# (everything except Fecha/Volumen) into the newly inserted row 85, then set
# the new row's own Fecha (D) and Volumen (J) values.
$ws.Range("A85").Value = $ws.Range("A86").Value2
$ws.Range("B85").Value = $ws.Range("B86").Value2
$ws.Range("C85").Value = $ws.Range("C86").Value2
$ws.Range("D85").Value = 45049
$ws.Range("E85").Value = $ws.Range("E86").Value2
$ws.Range("F85").Value = $ws.Range("F86").Value2
$ws.Range("G85").Value = $ws.Range("G86").Value2
$ws.Range("H85").Value = $ws.Range("H86").Value2
$ws.Range("I85").Value = $ws.Range("I86").Value2
$ws.Range("J85").Value = 35
$ws.Range("K85").Value = $ws.Range("K86").Value2
$ws.Range("L85").Value = $ws.Range("L86").Value2
$ws.Range("M85").Value = $ws.Range("M86").Value2
$ws.Range("N85").Value = $ws.Range("N86").Value2
$ws.Range("O85").Value = $ws.Range("O86").Value2
$ws.Range("P85").Value = $ws.Range("P86").Value2
$ws.Range("Q85").Value = $ws.Range("Q86").Value2
$ws.Range("R85").Value = $ws.Range("R86").Value2
